# Remove the obsolete data row (Equipment Number PCIU1438389 / Vessel "KOTA PERWIRA"
# / Voyage 0004E / WONumber 7032005141 / ReferenceNumber 7032005141 / BOLNumber
# DEL900004900), which was row 2 of Sheet1. Deleting the entire row shifts the
# remaining data rows up and lets Excel naturally prune the now-unused shared
# strings, matching the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
